# Auto-generated Excel COM-interop script
# Applies 2024-10-18 crime data update to violent-crime-full-year.xlsx
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 6158
$ws.Range("K3").Value = 6336
$ws.Range("C4").Value = 1851
$ws.Range("D4").Value = 1975
$ws.Range("K4").Value = 1331
$ws.Range("K6").Value = 6987
$ws.Range("C7").Value = 28396
$ws.Range("D7").Value = 28165
$ws.Range("K7").Value = 21262

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K2").Value = 71
$ws.Range("K7").Value = 270

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K6").Value = 471
$ws.Range("K7").Value = 1393

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 241
$ws.Range("K3").Value = 333
$ws.Range("K6").Value = 283
$ws.Range("K7").Value = 928

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 205
$ws.Range("K3").Value = 240
$ws.Range("K7").Value = 720

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K3").Value = 143
$ws.Range("K6").Value = 90
$ws.Range("K7").Value = 349

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("K2").Value = 23
$ws.Range("K7").Value = 79

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 187
$ws.Range("K7").Value = 622
$ws.Range("K8").Value = 1393
$ws.Range("K11").Value = 396
$ws.Range("K12").Value = 38
$ws.Range("K18").Value = 140
$ws.Range("K19").Value = 620
$ws.Range("K20").Value = 508
$ws.Range("K23").Value = 217
$ws.Range("K25").Value = 102
$ws.Range("K27").Value = 202
$ws.Range("K29").Value = 1159
$ws.Range("K30").Value = 79
$ws.Range("K33").Value = 928
$ws.Range("K34").Value = 121
$ws.Range("K35").Value = 33
$ws.Range("K36").Value = 275
$ws.Range("K37").Value = 720
$ws.Range("K40").Value = 48
$ws.Range("K42").Value = 787
$ws.Range("K43").Value = 179
$ws.Range("K49").Value = 116
$ws.Range("K52").Value = 564
$ws.Range("K53").Value = 270
$ws.Range("K54").Value = 413
$ws.Range("K57").Value = 80
$ws.Range("K60").Value = 126
$ws.Range("K62").Value = 7
$ws.Range("C63").Value = 280
$ws.Range("D63").Value = 354
$ws.Range("K63").Value = 62
$ws.Range("K64").Value = 133
$ws.Range("K67").Value = 836
$ws.Range("K68").Value = 58
$ws.Range("K70").Value = 36
$ws.Range("K72").Value = 106
$ws.Range("K73").Value = 185
$ws.Range("K76").Value = 288
$ws.Range("K78").Value = 238
$ws.Range("K79").Value = 533
$ws.Range("K80").Value = 74
$ws.Range("K84").Value = 166
$ws.Range("K85").Value = 985
$ws.Range("K89").Value = 315
$ws.Range("K91").Value = 243
$ws.Range("K93").Value = 81
$ws.Range("K94").Value = 287
$ws.Range("K97").Value = 167
$ws.Range("K98").Value = 104
$ws.Range("K99").Value = 349
$ws.Range("C101").Value = 28396
$ws.Range("D101").Value = 28165
$ws.Range("K101").Value = 21262

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K3").Value = 302
$ws.Range("K7").Value = 836

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K3").Value = 65
$ws.Range("K7").Value = 166

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("K6").Value = 59
$ws.Range("K7").Value = 116

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K6").Value = 225
$ws.Range("K7").Value = 413

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K3").Value = 417
$ws.Range("K6").Value = 327
$ws.Range("K7").Value = 1159

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K2").Value = 184
$ws.Range("K6").Value = 202
$ws.Range("K7").Value = 620

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K2").Value = 65
$ws.Range("K6").Value = 149
$ws.Range("K7").Value = 288

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K6").Value = 294
$ws.Range("K7").Value = 787

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K4").Value = 22
$ws.Range("K7").Value = 238

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("K3").Value = 76
$ws.Range("K7").Value = 217

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K3").Value = 117
$ws.Range("K7").Value = 243

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K2").Value = 177
$ws.Range("K6").Value = 135
$ws.Range("K7").Value = 533

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("K3").Value = 37
$ws.Range("K7").Value = 133

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 174
$ws.Range("K7").Value = 508

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("K4").Value = 17
$ws.Range("K7").Value = 140

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K3").Value = 82
$ws.Range("K6").Value = 62
$ws.Range("K7").Value = 275

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("K6").Value = 33
$ws.Range("K7").Value = 81

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 203
$ws.Range("K3").Value = 206
$ws.Range("K6").Value = 168
$ws.Range("K7").Value = 622

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("K3").Value = 32
$ws.Range("K7").Value = 121

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("K6").Value = 127
$ws.Range("K7").Value = 287

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("K2").Value = 36
$ws.Range("K7").Value = 102

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("K6").Value = 60
$ws.Range("K7").Value = 104

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K2").Value = 140
$ws.Range("K7").Value = 396

$ws = $wb.Worksheets.Item("Gold Coast")
$ws.Range("K6").Value = 20
$ws.Range("K7").Value = 33

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K2").Value = 62
$ws.Range("K7").Value = 185

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K2").Value = 55
$ws.Range("K3").Value = 49
$ws.Range("K7").Value = 187

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("K6").Value = 93
$ws.Range("K7").Value = 167

$ws = $wb.Worksheets.Item("O'Hare")
$ws.Range("K3").Value = 8
$ws.Range("K7").Value = 36

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K2").Value = 90
$ws.Range("K7").Value = 315

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K3").Value = 47
$ws.Range("K4").Value = 28
$ws.Range("K7").Value = 202

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("K2").Value = 25
$ws.Range("K7").Value = 58

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("K2").Value = 22
$ws.Range("K7").Value = 80

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("K4").Value = 10
$ws.Range("K7").Value = 126

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K4").Value = 24
$ws.Range("K6").Value = 72
$ws.Range("K7").Value = 179

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K3").Value = 341
$ws.Range("K6").Value = 240
$ws.Range("K7").Value = 985

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("K2").Value = 21
$ws.Range("K7").Value = 106

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("K4").Value = 7
$ws.Range("K7").Value = 74

$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Range("K3").Value = 21
$ws.Range("K7").Value = 48

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K3").Value = 163
$ws.Range("K6").Value = 201
$ws.Range("K7").Value = 564

$ws = $wb.Worksheets.Item("Beverly")
$ws.Range("K4").Value = 4
$ws.Range("K7").Value = 38

$ws = $wb.Worksheets.Item("Museum Campus")
$ws.Range("K6").Value = 5
$ws.Range("K7").Value = 7
